# Updated cryptos list on Sun Mar 31 09:26:11 UTC 2024 with GitHub Actions
# Applies per-cell value updates (prices, 1h volume %) and fixes the
# Chainlink / WrappedBTC row ordering (rows 17-18 swapped back).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.113.64"
$ws.Range("E2").Value = "  +0.22%  "

# Row 3
$ws.Range("D3").Value = "3.605.78"
$ws.Range("E3").Value = "  +2.84%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.98"
$ws.Range("E5").Value = "  +0.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.02"
$ws.Range("E6").Value = "  +0.29%  "

# Row 7
$ws.Range("E7").Value = "  +0.38%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  -0.75%  "

# Row 10
$ws.Range("E10").Value = "  -0.51%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.83"
$ws.Range("E11").Value = "  -0.32%  "

# Row 12
$ws.Range("E12").Value = "  +1.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.56"
$ws.Range("E13").Value = "  +0.24%  "

# Row 14
$ws.Range("D14").Value = "4.185.24"
$ws.Range("E14").Value = "  +3.17%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.13"
$ws.Range("E15").Value = "  +4.19%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "592.32"
$ws.Range("E16").Value = "  -2.15%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "70.264.24"
$ws.Range("E17").Value = "  +0.36%  "

# Row 18
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.18"
$ws.Range("E18").Value = "  +1.08%  "

# Row 19
$ws.Range("D19").Value = "3.608.64"
$ws.Range("E19").Value = "  +2.95%  "

# Row 21
$ws.Range("E21").Value = "  +0.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.70"
$ws.Range("E22").Value = "  -2.37%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.15"
$ws.Range("E23").Value = "  +0.36%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.06"
$ws.Range("E24").Value = "  -2.27%  "

# Row 25
$ws.Range("E25").Value = "  +0.49%  "

# Row 26
$ws.Range("E26").Value = "  -0.95%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.81"
$ws.Range("E27").Value = "  -1.17%  "

# Row 28
$ws.Range("E28").Value = "  -1.32%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.02"
$ws.Range("E29").Value = "  +1.36%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.80"
$ws.Range("E30").Value = "  +4.41%  "

# Row 31
$ws.Range("E31").Value = "  +0.44%  "

# Row 32
$ws.Range("E32").Value = "  -2.63%  "

# Row 33
$ws.Range("E33").Value = "  +1.65%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.27"
$ws.Range("E34").Value = "  -0.28%  "

# Row 35
$ws.Range("E35").Value = "  +7.81%  "

# Row 36
$ws.Range("D36").Value = "3.926.37"
$ws.Range("E36").Value = "  +5.27%  "

# Row 37
$ws.Range("E37").Value = "  +3.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "529.95"
$ws.Range("E38").Value = "  +6.35%  "

# Row 39
$ws.Range("E39").Value = "  -0.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.12"
$ws.Range("E40").Value = "  +1.00%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.391"
$ws.Range("E41").Value = "  -0.09%  "

# Row 42
$ws.Range("E42").Value = "  -1.11%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.134"
$ws.Range("E43").Value = "  -1.39%  "

# Row 44
$ws.Range("E44").Value = "  -0.50%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.39"
$ws.Range("E45").Value = "  +2.07%  "

# Row 46
$ws.Range("E46").Value = "  +1.22%  "

# Row 47
$ws.Range("E47").Value = "  +0.70%  "

# Row 48
$ws.Range("E48").Value = "  -1.38%  "

# Row 49
$ws.Range("E49").Value = "  -0.11%  "

# Row 50
$ws.Range("E50").Value = "  +4.11%  "

# Row 51
$ws.Range("E51").Value = "  +3.77%  "
